$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Columns A..G already exist; update their text to the new column names.
# Columns H..M are brand-new; set values and copy the header formatting (bold, centered, thin border).
$ws.Range("A1").Value = "index"

$ws.Range("B1").Value = "Temperature"

$ws.Range("C1").Value = "MashTime"

$ws.Range("D1").Value = "SolidosFermentaveis"

$ws.Range("E1").Value = "SolidosNaoFermentaveis"

$ws.Range("F1").Value = "SolidosTotais"

$ws.Range("G1").Value = "PercFermentaveis"

$ws.Range("H1").Value = "Extrato"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1

$ws.Range("I1").Value = "MashingEfficiency"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4160
$ws.Range("I1").Borders.LineStyle = 1

$ws.Range("J1").Value = "Dp1"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").VerticalAlignment = -4160
$ws.Range("J1").Borders.LineStyle = 1

$ws.Range("K1").Value = "Dp2"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("K1").VerticalAlignment = -4160
$ws.Range("K1").Borders.LineStyle = 1

$ws.Range("L1").Value = "Dp3"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").HorizontalAlignment = -4108
$ws.Range("L1").VerticalAlignment = -4160
$ws.Range("L1").Borders.LineStyle = 1

$ws.Range("M1").Value = "Dp4Plus"
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").HorizontalAlignment = -4108
$ws.Range("M1").VerticalAlignment = -4160
$ws.Range("M1").Borders.LineStyle = 1

# --- Data rows 2..9 ---
# Column A keeps the same stat-label text (count/mean/std/min/25%/50%/75%/max).
# Columns B..M get the new numeric values for the new dataset.
$ws.Range("A2").Value = "count"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 1000

$ws.Range("A3").Value = "mean"
$ws.Range("B3").Value = 65
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 96.9409299426404
$ws.Range("E3").Value = 48.49040226979136
$ws.Range("F3").Value = 152.4615906654977
$ws.Range("G3").Value = 63.47926899177205
$ws.Range("H3").Value = 14.54313322124318
$ws.Range("I3").Value = 95.29886035261634
$ws.Range("J3").Value = 7.913637169029514
$ws.Range("K3").Value = 45.49775780026358
$ws.Range("L3").Value = 10.06787402247895
$ws.Range("M3").Value = 31.81959136084431

$ws.Range("A4").Value = "std"
$ws.Range("B4").Value = 9.998494677726423
$ws.Range("C4").Value = 29.99548403317927
$ws.Range("D4").Value = 22.09488519857674
$ws.Range("E4").Value = 13.31724655041774
$ws.Range("F4").Value = 1.123328699128405
$ws.Range("G4").Value = 14.22205674278546
$ws.Range("H4").Value = 2.275842876167337
$ws.Range("I4").Value = 14.6476119442378
$ws.Range("J4").Value = 1.250674003895388
$ws.Range("K4").Value = 12.80693754054486
$ws.Range("L4").Value = 2.075530562340822
$ws.Range("M4").Value = 8.867311140844556

$ws.Range("A5").Value = "min"
$ws.Range("B5").Value = 32.09473268508106
$ws.Range("C5").Value = 1.284198055243166
$ws.Range("D5").Value = 15.50705457086778
$ws.Range("E5").Value = 13.59428172819732
$ws.Range("F5").Value = 148.10697344438
$ws.Range("G5").Value = 10.4701718023367
$ws.Range("H5").Value = 3.600697344437997
$ws.Range("I5").Value = 24.31146394190683
$ws.Range("J5").Value = 3.46254299624347
$ws.Range("K5").Value = 7.00762880609323
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 9.150886220158162

$ws.Range("A6").Value = "25%"
$ws.Range("B6").Value = 58.26296131778599
$ws.Range("C6").Value = 79.78888395335795
$ws.Range("D6").Value = 89.73933244461371
$ws.Range("E6").Value = 40.11173855912033
$ws.Range("F6").Value = 152.163736052687
$ws.Range("G6").Value = 58.98729052204998
$ws.Range("H6").Value = 15.15110282411983
$ws.Range("I6").Value = 99.64800337717145
$ws.Range("J6").Value = 7.199932893177007
$ws.Range("K6").Value = 39.58943849195545
$ws.Range("L6").Value = 10.6950767294104
$ws.Range("M6").Value = 26.20262024773336

$ws.Range("A7").Value = "50%"
$ws.Range("B7").Value = 65
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = 105.8768879966164
$ws.Range("E7").Value = 44.86397495528115
$ws.Range("F7").Value = 152.9140553288643
$ws.Range("G7").Value = 69.23879595889238
$ws.Range("H7").Value = 15.28698364340404
$ws.Range("I7").Value = 99.99916469286694
$ws.Range("J7").Value = 7.670735716254943
$ws.Range("K7").Value = 50.65271883845252
$ws.Range("L7").Value = 10.70235330613033
$ws.Range("M7").Value = 29.34608123116248

$ws.Range("A8").Value = "75%"
$ws.Range("B8").Value = 71.73703868221402
$ws.Range("C8").Value = 120.211116046642
$ws.Range("D8").Value = 112.2942878222007
$ws.Range("E8").Value = 51.42766546338035
$ws.Range("F8").Value = 153.2103440366889
$ws.Range("G8").Value = 73.29419466744451
$ws.Range("H8").Value = 15.32053653781785
$ws.Range("I8").Value = 99.99999586186635
$ws.Range("J8").Value = 8.316078059309714
$ws.Range("K8").Value = 55.29412598582952
$ws.Range("L8").Value = 10.72136415644584
$ws.Range("M8").Value = 33.72328583140195

$ws.Range("A9").Value = "max"
$ws.Range("B9").Value = 97.90526731491926
$ws.Range("C9").Value = 198.7158019447578
$ws.Range("D9").Value = 114.5695897190051
$ws.Range("E9").Value = 85.74294277731047
$ws.Range("F9").Value = 153.3148027818198
$ws.Range("G9").Value = 74.7283286676809
$ws.Range("H9").Value = 15.33148027346516
$ws.Range("I9").Value = 99.99999999996247
$ws.Range("J9").Value = 10.74933242774137
$ws.Range("K9").Value = 56.95111272001579
$ws.Range("L9").Value = 10.79451731378597
$ws.Range("M9").Value = 57.03254907861945
